$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 409, shifting existing rows 409:471 down to 410:472
$ws.Rows.Item(409).Insert()

# Populate the newly inserted row 409 with the new data record
$ws.Cells.Item(409, 1).Value = 10
$ws.Cells.Item(409, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(409, 3).Value = "La Araucanía"
$ws.Cells.Item(409, 4).Value = 45034
$ws.Cells.Item(409, 5).Value = 9
$ws.Cells.Item(409, 6).Value = 100112044
$ws.Cells.Item(409, 7).Value = "Perejil"
$ws.Cells.Item(409, 8).Value = "Sin especificar"
$ws.Cells.Item(409, 9).Value = "Primera"
$ws.Cells.Item(409, 10).Value = 65
$ws.Cells.Item(409, 11).Value = 4000
$ws.Cells.Item(409, 12).Value = 4000
$ws.Cells.Item(409, 13).Value = 4000
$ws.Cells.Item(409, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(409, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(409, 16).Value = 1333
$ws.Cells.Item(409, 17).Value = 3
$ws.Cells.Item(409, 18).Value = "Hortaliza"
